$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.966.16"
Set-TextValue "E2" "  -2.88%  "

Set-TextValue "D3" "1.743.88"
Set-TextValue "E3" "  -0.87%  "

Set-TextValue "E4" "  -0.29%  "

Set-TextValue "D6" "1.0000"
Set-TextValue "E6" "  -0.20%  "

Set-TextValue "D7" "0.4992"
Set-TextValue "E7" "  +3.24%  "

Set-TextValue "D8" "0.3581"
Set-TextValue "E8" "  +0.64%  "

Set-TextValue "D9" "42.63"
Set-TextValue "E9" "  -1.45%  "

Set-TextValue "D10" "0.07268"
Set-TextValue "E10" "  -3.34%  "

Set-TextValue "E11" "  -1.79%  "

Set-TextValue "E12" "  -0.23%  "

Set-TextValue "D13" "20.08"
Set-TextValue "E13" "  -2.26%  "

Set-TextValue "D14" "5.983"
Set-TextValue "E14" "  -1.98%  "

Set-TextValue "D15" "1.741.86"
Set-TextValue "E15" "  -0.97%  "

Set-TextValue "D16" "6.858"
Set-TextValue "E16" "  -3.64%  "

Set-TextValue "D17" "86.59"
Set-TextValue "E17" "  -7.19%  "

Set-TextValue "D18" "0.00001036"
Set-TextValue "E18" "  -4.22%  "

Set-TextValue "D19" "0.06397"
Set-TextValue "E19" "  -0.55%  "

Set-TextValue "D20" "0.9987"
Set-TextValue "E20" "  -0.25%  "

Set-TextValue "E21" "  -1.16%  "

Set-TextValue "D22" "5.738"
Set-TextValue "E22" "  -1.03%  "

Set-TextValue "D23" "27.010.89"
Set-TextValue "E23" "  -2.88%  "

Set-TextValue "D24" "11.31"
Set-TextValue "E24" "  +1.90%  "

Set-TextValue "D25" "2.047"
Set-TextValue "E25" "  -5.16%  "

Set-TextValue "D26" "153.72"
Set-TextValue "E26" "  -5.94%  "

Set-TextValue "D27" "19.93"
Set-TextValue "E27" "  -1.08%  "

Set-TextValue "D28" "1.938.43"
Set-TextValue "E28" "  -1.14%  "

Set-TextValue "D29" "2.145"
Set-TextValue "E29" "  -2.68%  "

Set-TextValue "D30" "120.80"
Set-TextValue "E30" "  -1.55%  "

Set-TextValue "D31" "1.064"
Set-TextValue "E31" "  +0.77%  "

Set-TextValue "D32" "0.09462"
Set-TextValue "E32" "  +0.30%  "

Set-TextValue "D33" "3.566"
Set-TextValue "E33" "  -2.48%  "

Set-TextValue "D34" "5.395"
Set-TextValue "E34" "  -2.80%  "

Set-TextValue "D35" "0.02204"
Set-TextValue "E35" "  -2.73%  "

Set-TextValue "D36" "0.05916"
Set-TextValue "E36" "  -1.02%  "

Set-TextValue "D37" "11.12"
Set-TextValue "E37" "  -4.25%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D38" "0.2000"
Set-TextValue "E38" "  -3.10%  "

$ws.Range("B39").Value = "WEMIXTOKEN"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D39" "1.422"
Set-TextValue "E39" "  -0.91%  "

Set-TextValue "D40" "4.764"
Set-TextValue "E40" "  -2.44%  "

Set-TextValue "D41" "0.6031"
Set-TextValue "E41" "  -1.83%  "

Set-TextValue "E42" "  -0.23%  "

Set-TextValue "D43" "1.111"
Set-TextValue "E43" "  -5.83%  "

Set-TextValue "D44" "7.472"
Set-TextValue "E44" "  -3.54%  "

Set-TextValue "D45" "12.80"
Set-TextValue "E45" "  -2.13%  "

Set-TextValue "D46" "3.582"
Set-TextValue "E46" "  -4.04%  "

Set-TextValue "D47" "0.5654"
Set-TextValue "E47" "  -2.37%  "

Set-TextValue "D48" "119.99"

Set-TextValue "D49" "1.863"
Set-TextValue "E49" "  -3.27%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D50" "0.06676"
Set-TextValue "E50" "  -1.65%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue "D51" "1.101"
Set-TextValue "E51" "  -4.22%  "
